$d = $word.ActiveDocument
$t = $d.Tables(1)
$pairs = @(
    @("28+14=", "61+6="),
    @("75-28=", "88-74="),
    @("70-66=", "38+31="),
    @("33-0=", "60+25="),
    @("49+32=", "38+3="),
    @("23+49=", "7+18="),
    @("37+4=", "44-24="),
    @("94-91=", "87+9="),
    @("24+64=", "42+11="),
    @("44-37=", "69-40="),
    @("67-56=", "87-63="),
    @("40+19=", "50-39="),
    @("63-32=", "46+47="),
    @("17-4=", "70-29="),
    @("8+0=", "0+24="),
    @("59-57=", "40-27="),
    @("18+73=", "79-57="),
    @("26+11=", "79-44="),
    @("15+21=", "27+50="),
    @("58+9=", "5+74="),
    @("23-0=", "93-65="),
    @("58-24=", "19+70="),
    @("72-20=", "31+38="),
    @("11+57=", "56+19="),
    @("36+41=", "67-16="),
    @("87-73=", "0+17="),
    @("14+6=", "69+16="),
    @("95-65=", "28+2="),
    @("26+38=", "28+27="),
    @("48-1=", "8+67="),
    @("29-1=", "38+13="),
    @("61+12=", "7+0="),
    @("14-6=", "30+51="),
    @("69-18=", "94-1="),
    @("96-37=", "73+13="),
    @("63-45=", "91-61="),
    @("74+1=", "27+41="),
    @("5+51=", "5+60="),
    @("42+49=", "77-62="),
    @("36+52=", "62+11="),
    @("66-5=", "53-4="),
    @("65+29=", "64-23="),
    @("88-66=", "10+56="),
    @("18+73=", "77-15="),
    @("57-24=", "6+70="),
    @("47+43=", "18-4="),
    @("57-23=", "2+34="),
    @("0+27=", "9+52="),
    @("2+55=", "26+8="),
    @("46-37=", "50+33="),
    @("80-19=", "54+6="),
    @("78-41=", "81+15="),
    @("54-22=", "44+12="),
    @("88-36=", "33-28="),
    @("35+46=", "36-4="),
    @("14+32=", "2+48="),
    @("63-7=", "70+4="),
    @("45+26=", "48+47="),
    @("36+22=", "22+68="),
    @("29+68=", "15+8="),
    @("13+77=", "42+2="),
    @("59-51=", "80-64="),
    @("32+42=", "12+82="),
    @("18+81=", "78-63="),
    @("98-77=", "14+7="),
    @("51+15=", "18+31="),
    @("66-14=", "67+5="),
    @("5+0=", "98-92="),
    @("45-19=", "45-6="),
    @("36-1=", "21+7="),
    @("10-1=", "64+22="),
    @("37-29=", "2+93="),
    @("47+5=", "52+47="),
    @("92-83=", "82-74="),
    @("59-24=", "45-20="),
    @("24+8=", "95-1="),
    @("48+35=", "75+14="),
    @("27+64=", "6-4="),
    @("92-18=", "58-43="),
    @("15+27=", "67-50="),
    @("0+72=", "16+16="),
    @("59+20=", "10+84="),
    @("5+26=", "41+10="),
    @("98-52=", "15+2="),
    @("30+27=", "25-10="),
    @("86-10=", "60+29="),
    @("27-22=", "45-15="),
    @("30+20=", "87-76="),
    @("59-17=", "95-28="),
    @("15+47=", "81-78="),
    @("63+9=", "14+4="),
    @("56-41=", "65+0="),
    @("43-12=", "49-43="),
    @("1+14=", "44-8="),
    @("83-59=", "49+18="),
    @("11+67=", "3+41="),
    @("3+73=", "97-58="),
    @("57-57=", "30+5="),
    @("42+32=", "5+39="),
    @("80-23=", "26+38=")
)

$cols = $t.Columns.Count
$applied = 0
$skipped = 0
for ($i = 0; $i -lt $pairs.Count; $i++) {
    $row = [int][math]::Floor($i / $cols) + 1
    $col = ($i % $cols) + 1
    $old = $pairs[$i][0]
    $new = $pairs[$i][1]
    $cell = $t.Cell($row, $col)
    $rng = $cell.Range
    # Replace=1 (wdReplaceOne) + Wrap=0 (wdFindStop) keeps the Find strictly
    # scoped to this cell's Range, which matters because some expressions
    # (e.g. "18+73=") repeat at different positions with different targets.
    $res = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 0, $false, $new, 1)
    if ($res) {
        $applied = $applied + 1
    } else {
        $skipped = $skipped + 1
        Write-Host "WARN: cell($row,$col) expected '$old' but replace failed; current text:" $cell.Range.Text
    }
}
Write-Host "Replacements applied:" $applied "skipped:" $skipped
